$d = $word.ActiveDocument

# 1. "differences in disease in either experiment due to irrigation regime, N rates or
#    the interaction of the two treatments in either season. This suggests"
#    -> "difference in the incidence of tiller sheath blight due to irrigation, tiller
#    and leaf sheath blight did differ significantly by irrigation treament but leaf
#    sheath blight severity did not. Our findings suggests"
$d.Content.Find.Execute(
    "differences in disease in either experiment due to irrigation regime, N rates or the interaction of the two treatments in either season. This suggests",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "difference in the incidence of tiller sheath blight due to irrigation, tiller and leaf sheath blight did differ significantly by irrigation treament but leaf sheath blight severity did not. Our findings suggests",
    2)

# 2. Update the generated timestamp in the colophon paragraph.
$d.Content.Find.Execute(
    "This report was generated on 2018-04-12 15:59:46 using the following computational environment and dependencies:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This report was generated on 2018-04-12 16:08:53 using the following computational environment and dependencies:",
    2)

# 3. Update the git HEAD commit hash/message shown in the colophon.
$d.Content.Find.Execute(
    "## Head:     [ee4a88b] 2018-04-12: Fix references",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## Head:     [b7fc71e] 2018-04-12: Update paper and add Word output",
    2)
